$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4691.8125
$ws.Range("J74").Value = 4994.25
$ws.Range("L74").Value = 4994.25
$ws.Range("N74").Value = -6866.25

$ws.Range("H76").Value = 5615.8125
$ws.Range("I76").Value = 3936.3333
$ws.Range("J76").Value = 6003.385
$ws.Range("K76").Value = 3936.3333
$ws.Range("L76").Value = 6003.385
$ws.Range("M76").Value = -3621.3333
$ws.Range("N76").Value = -6633.385

$ws.Range("H77").Value = 4691.8125
$ws.Range("J77").Value = 4994.25
$ws.Range("L77").Value = 24971.25
$ws.Range("N77").Value = -34331.25

$ws.Range("H79").Value = 5615.8125
$ws.Range("I79").Value = 3936.3333
$ws.Range("J79").Value = 6003.385
$ws.Range("K79").Value = 3936.3333
$ws.Range("L79").Value = 6003.385
$ws.Range("M79").Value = -2844.3333
$ws.Range("N79").Value = -8187.385

$ws.Range("H138").Value = 262028.22
$ws.Range("I138").Value = 694477.0600000001
$ws.Range("J138").Value = 5261.7188
$ws.Range("K138").Value = 2083431.18
$ws.Range("L138").Value = 15785.1564
$ws.Range("M138").Value = -2078291.18
$ws.Range("N138").Value = -26065.1564

$ws.Range("H139").Value = 137466.8
$ws.Range("J139").Value = 137466.8
$ws.Range("L139").Value = 137466.8
$ws.Range("N139").Value = -147746.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2912.6052
$ws.Range("J2").Value = 3182.182
$ws.Range("L2").Value = 3182.182
$ws.Range("N2").Value = -3408.182

$ws.Range("H116").Value = 2912.6052
$ws.Range("J116").Value = 3182.182
$ws.Range("L116").Value = 3182.182
$ws.Range("N116").Value = -7770.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2912.6052
$ws.Range("J3").Value = 3182.182
$ws.Range("L3").Value = 3182.182
$ws.Range("N3").Value = -3410.182

$ws.Range("H86").Value = 6581.15
$ws.Range("I86").Value = 7798.0713
$ws.Range("K86").Value = 7798.0713
$ws.Range("M86").Value = -6675.0713

$ws.Range("H89").Value = 6581.15
$ws.Range("I89").Value = 7798.0713
$ws.Range("K89").Value = 38990.35649999999
$ws.Range("M89").Value = -33374.35649999999

$ws.Range("H94").Value = 3030.9707
$ws.Range("I94").Value = 2222.6785
$ws.Range("J94").Value = 6803
$ws.Range("K94").Value = 2222.6785
$ws.Range("L94").Value = 6803
$ws.Range("M94").Value = -1771.6785
$ws.Range("N94").Value = -7705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2289
$ws.Range("I94").Value = 1002
$ws.Range("J94").Value = 2432
$ws.Range("K94").Value = 1002
$ws.Range("L94").Value = 2432
$ws.Range("M94").Value = -551
$ws.Range("N94").Value = -3334

$ws.Range("H134").Value = 1872.8959
$ws.Range("I134").Value = 1453
$ws.Range("K134").Value = 4359
$ws.Range("M134").Value = -1824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3386.6365
$ws.Range("I46").Value = 805.5
$ws.Range("J46").Value = 5537.5835
$ws.Range("K46").Value = 2416.5
$ws.Range("L46").Value = 16612.7505
$ws.Range("M46").Value = -2325.5
$ws.Range("N46").Value = -16794.7505

$ws.Range("H64").Value = 2685.4285
$ws.Range("J64").Value = 3899.5
$ws.Range("L64").Value = 11698.5
$ws.Range("N64").Value = -12238.5

$ws.Range("H67").Value = 2685.4285
$ws.Range("J67").Value = 3899.5
$ws.Range("L67").Value = 11698.5
$ws.Range("N67").Value = -13570.5

$ws.Range("H82").Value = 7370
$ws.Range("J82").Value = 10555
$ws.Range("L82").Value = 31665
$ws.Range("N82").Value = -32477

$ws.Range("H85").Value = 7370
$ws.Range("J85").Value = 10555
$ws.Range("L85").Value = 31665
$ws.Range("N85").Value = -34473

$ws.Range("H97").Value = 39614.125
$ws.Range("I97").Value = 43823.285
$ws.Range("K97").Value = 131469.855
$ws.Range("M97").Value = -130973.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12288.941
$ws.Range("I80").Value = 18171.777
$ws.Range("J80").Value = 5670.75
$ws.Range("K80").Value = 18171.777
$ws.Range("L80").Value = 5670.75
$ws.Range("M80").Value = -17173.777
$ws.Range("N80").Value = -7666.75

$ws.Range("H83").Value = 12288.941
$ws.Range("I83").Value = 18171.777
$ws.Range("J83").Value = 5670.75
$ws.Range("K83").Value = 90858.88499999999
$ws.Range("L83").Value = 28353.75
$ws.Range("M83").Value = -85866.88499999999
$ws.Range("N83").Value = -38337.75

$ws.Range("H122").Value = 7733.794
$ws.Range("I122").Value = 5350.174
$ws.Range("K122").Value = 16050.522
$ws.Range("M122").Value = -13600.522

$ws.Range("H126").Value = 17867
$ws.Range("I126").Value = 22489.5
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 67468.5
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -64998.5
$ws.Range("N126").Value = -16938.5

$ws.Range("H133").Value = 69998
$ws.Range("J133").Value = 69998
$ws.Range("L133").Value = 69998
$ws.Range("N133").Value = -80118

$ws.Range("H135").Value = 169555
$ws.Range("J135").Value = 169555
$ws.Range("L135").Value = 169555
$ws.Range("N135").Value = -179695

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 44532.75
$ws.Range("I61").Value = 2165.5
$ws.Range("K61").Value = 2165.5
$ws.Range("M61").Value = -1963.5

$ws.Range("H82").Value = 1816.7368
$ws.Range("J82").Value = 2033
$ws.Range("L82").Value = 2033
$ws.Range("N82").Value = -2755

$ws.Range("H85").Value = 1816.7368
$ws.Range("J85").Value = 2033
$ws.Range("L85").Value = 2033
$ws.Range("N85").Value = -4529

$ws.Range("H92").Value = 59999.5
$ws.Range("J92").Value = 59999.5
$ws.Range("L92").Value = 59999.5
$ws.Range("N92").Value = -64991.5

$ws.Range("H93").Value = 4207.3184
$ws.Range("I93").Value = 4435
$ws.Range("J93").Value = 2765.3333
$ws.Range("K93").Value = 4435
$ws.Range("L93").Value = 2765.3333
$ws.Range("M93").Value = -3187
$ws.Range("N93").Value = -5261.3333

$ws.Range("H113").Value = 44532.75
$ws.Range("I113").Value = 2165.5
$ws.Range("K113").Value = 2165.5
$ws.Range("M113").Value = 4.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 433222.34
$ws.Range("I62").Value = 639833.2
$ws.Range("K62").Value = 639833.2
$ws.Range("M62").Value = -639209.2

$ws.Range("H65").Value = 433222.34
$ws.Range("I65").Value = 639833.2
$ws.Range("K65").Value = 3199166
$ws.Range("M65").Value = -3196046

$ws.Range("H122").Value = 5259.7
$ws.Range("I122").Value = 3143.25
$ws.Range("K122").Value = 9429.75
$ws.Range("M122").Value = -6979.75

$ws.Range("H136").Value = 246847.33
$ws.Range("I136").Value = 263213.25
$ws.Range("J136").Value = 5450
$ws.Range("K136").Value = 789639.75
$ws.Range("L136").Value = 16350
$ws.Range("M136").Value = -787089.75
$ws.Range("N136").Value = -21450
